# Generate Report for Handoff
# ----------------------------------------------------------------------------
# The report's three data rows (one per source file) get re-sorted by latest
# handoff date, and the "fc53b026..." file - which has just been hand off
# again - moves to the bottom of each table with a refreshed status and
# handoff timestamps.
#
# Row order (old -> new):
#   old: fc53b026..., ffff5b016cb9..., ffffff452bf4df...
#   new: ffff5b016cb9..., ffffff452bf4df..., fc53b026...
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: match on exact single-cell address string, e.g. "$A$2"
function Set-Cell {
    param($ws, [string]$addr, [string]$text)
    $ws.Range($addr).Value = $text
}

function Set-HyperlinkDisplay {
    param($ws, [string]$dollarAddr, [string]$text)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $dollarAddr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-Cell $wsOverview "A2" "ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md"
Set-Cell $wsOverview "B2" "Handed back: in sync with en-US"
Set-Cell $wsOverview "C2" "Handed back: in sync with en-US"
Set-Cell $wsOverview "D2" "2016-03-21 03:42:16"

Set-Cell $wsOverview "A3" "ffffff452bf4df-9ed6-4c2a-b71c-4d9845084bfc.md"
Set-Cell $wsOverview "B3" "Handed back: in sync with en-US"
Set-Cell $wsOverview "C3" "Handed back: in sync with en-US"
Set-Cell $wsOverview "D3" "2016-03-21 03:42:16"

Set-Cell $wsOverview "A4" "fc53b026-72f9-4f60-980f-e4271f141c78.md"
Set-Cell $wsOverview "B4" "Ready for handoff"
Set-Cell $wsOverview "C4" "Ready for handoff"
Set-Cell $wsOverview "D4" "2016-03-21 03:45:47"

Set-HyperlinkDisplay $wsOverview '$A$2' "ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md"
Set-HyperlinkDisplay $wsOverview '$A$3' "ffffff452bf4df-9ed6-4c2a-b71c-4d9845084bfc.md"
Set-HyperlinkDisplay $wsOverview '$A$4' "fc53b026-72f9-4f60-980f-e4271f141c78.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-Cell $wsZh "A2" "ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md"
Set-Cell $wsZh "B2" ".md"
Set-Cell $wsZh "C2" "Handed back: in sync with en-US"
Set-Cell $wsZh "D2" "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"
Set-Cell $wsZh "E2" "2016-03-21 03:42:07"
Set-Cell $wsZh "F2" "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
Set-Cell $wsZh "G2" "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"
Set-Cell $wsZh "H2" "2016-03-21 03:42:47"
Set-Cell $wsZh "J2" "Include"

Set-Cell $wsZh "A3" "ffffff452bf4df-9ed6-4c2a-b71c-4d9845084bfc.md"
Set-Cell $wsZh "B3" ".md"
Set-Cell $wsZh "C3" "Handed back: in sync with en-US"
Set-Cell $wsZh "D3" "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"
Set-Cell $wsZh "E3" "2016-03-21 03:42:07"
Set-Cell $wsZh "F3" "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
Set-Cell $wsZh "G3" "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"
Set-Cell $wsZh "H3" "2016-03-21 03:42:47"
Set-Cell $wsZh "J3" "Include"

Set-Cell $wsZh "A4" "fc53b026-72f9-4f60-980f-e4271f141c78.md"
Set-Cell $wsZh "B4" ".md"
Set-Cell $wsZh "C4" "Ready for handoff"
Set-Cell $wsZh "D4" "fc53b026-72f9-4f60-980f-e4271f141c78.c7d870cde34c605621f64a8f4e47b678a57047a0.zh-cn.xlf"
Set-Cell $wsZh "E4" "2016-03-21 03:45:39"
Set-Cell $wsZh "F4" "fc53b026-72f9-4f60-980f-e4271f141c78.md"
Set-Cell $wsZh "G4" "fc53b026-72f9-4f60-980f-e4271f141c78.c7d870cde34c605621f64a8f4e47b678a57047a0.zh-cn.xlf"
Set-Cell $wsZh "H4" "2016-03-21 03:44:45"
Set-Cell $wsZh "J4" "Include"

Set-HyperlinkDisplay $wsZh '$A$2' "ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md"
Set-HyperlinkDisplay $wsZh '$D$2' "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"
Set-HyperlinkDisplay $wsZh '$F$2' "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
Set-HyperlinkDisplay $wsZh '$G$2' "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"

Set-HyperlinkDisplay $wsZh '$A$3' "ffffff452bf4df-9ed6-4c2a-b71c-4d9845084bfc.md"
Set-HyperlinkDisplay $wsZh '$D$3' "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"
Set-HyperlinkDisplay $wsZh '$F$3' "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
Set-HyperlinkDisplay $wsZh '$G$3' "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.zh-cn.xlf"

Set-HyperlinkDisplay $wsZh '$A$4' "fc53b026-72f9-4f60-980f-e4271f141c78.md"
Set-HyperlinkDisplay $wsZh '$D$4' "fc53b026-72f9-4f60-980f-e4271f141c78.c7d870cde34c605621f64a8f4e47b678a57047a0.zh-cn.xlf"
Set-HyperlinkDisplay $wsZh '$F$4' "fc53b026-72f9-4f60-980f-e4271f141c78.md"
Set-HyperlinkDisplay $wsZh '$G$4' "fc53b026-72f9-4f60-980f-e4271f141c78.c7d870cde34c605621f64a8f4e47b678a57047a0.zh-cn.xlf"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-Cell $wsDe "A2" "ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md"
Set-Cell $wsDe "B2" ".md"
Set-Cell $wsDe "C2" "Handed back: in sync with en-US"
Set-Cell $wsDe "D2" "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"
Set-Cell $wsDe "E2" "2016-03-21 03:42:16"
Set-Cell $wsDe "F2" "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
Set-Cell $wsDe "G2" "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"
Set-Cell $wsDe "H2" "2016-03-21 03:43:00"
Set-Cell $wsDe "J2" "Include"

Set-Cell $wsDe "A3" "ffffff452bf4df-9ed6-4c2a-b71c-4d9845084bfc.md"
Set-Cell $wsDe "B3" ".md"
Set-Cell $wsDe "C3" "Handed back: in sync with en-US"
Set-Cell $wsDe "D3" "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"
Set-Cell $wsDe "E3" "2016-03-21 03:42:16"
Set-Cell $wsDe "F3" "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
Set-Cell $wsDe "G3" "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"
Set-Cell $wsDe "H3" "2016-03-21 03:43:00"
Set-Cell $wsDe "J3" "Include"

Set-Cell $wsDe "A4" "fc53b026-72f9-4f60-980f-e4271f141c78.md"
Set-Cell $wsDe "B4" ".md"
Set-Cell $wsDe "C4" "Ready for handoff"
Set-Cell $wsDe "D4" "fc53b026-72f9-4f60-980f-e4271f141c78.c7d870cde34c605621f64a8f4e47b678a57047a0.de-de.xlf"
Set-Cell $wsDe "E4" "2016-03-21 03:45:47"
Set-Cell $wsDe "F4" "fc53b026-72f9-4f60-980f-e4271f141c78.md"
Set-Cell $wsDe "G4" "fc53b026-72f9-4f60-980f-e4271f141c78.c7d870cde34c605621f64a8f4e47b678a57047a0.de-de.xlf"
Set-Cell $wsDe "H4" "2016-03-21 03:44:58"
Set-Cell $wsDe "J4" "Include"

Set-HyperlinkDisplay $wsDe '$A$2' "ffff5b016cb9-d55a-4ed9-a3d0-781c78144e0b.md"
Set-HyperlinkDisplay $wsDe '$D$2' "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"
Set-HyperlinkDisplay $wsDe '$F$2' "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
Set-HyperlinkDisplay $wsDe '$G$2' "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"

Set-HyperlinkDisplay $wsDe '$A$3' "ffffff452bf4df-9ed6-4c2a-b71c-4d9845084bfc.md"
Set-HyperlinkDisplay $wsDe '$D$3' "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"
Set-HyperlinkDisplay $wsDe '$F$3' "144dc3ed-8811-48a6-98eb-0053a6c9080a.md"
Set-HyperlinkDisplay $wsDe '$G$3' "144dc3ed-8811-48a6-98eb-0053a6c9080a.672ae3b566878e1a2e66c7d4276d1cb4ddeb1200.de-de.xlf"

Set-HyperlinkDisplay $wsDe '$A$4' "fc53b026-72f9-4f60-980f-e4271f141c78.md"
Set-HyperlinkDisplay $wsDe '$D$4' "fc53b026-72f9-4f60-980f-e4271f141c78.c7d870cde34c605621f64a8f4e47b678a57047a0.de-de.xlf"
Set-HyperlinkDisplay $wsDe '$F$4' "fc53b026-72f9-4f60-980f-e4271f141c78.md"
Set-HyperlinkDisplay $wsDe '$G$4' "fc53b026-72f9-4f60-980f-e4271f141c78.c7d870cde34c605621f64a8f4e47b678a57047a0.de-de.xlf"

Write-Output "Report regenerated for handoff."
